$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a numeric-looking string to a cell while preserving its
# original "text" storage (these Price cells are stored as inline strings,
# not numbers, so we must avoid Excel's automatic number coercion and any
# leftover formatting side-effects).
function Set-TextValue {
    param($Sheet, [string]$Addr, [string]$NewValue)

    $rng = $Sheet.Range($Addr)
    $originalStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = $originalStyle
}

Set-TextValue $ws "D2"  "243.04"
Set-TextValue $ws "D3"  "23.12"
Set-TextValue $ws "D4"  "5.419"
Set-TextValue $ws "D5"  "0.05892"
Set-TextValue $ws "D6"  "3.439"
Set-TextValue $ws "D7"  "6.545"
Set-TextValue $ws "D8"  "0.8106"
Set-TextValue $ws "D9"  "0.9379"
Set-TextValue $ws "D11" "0.07391"
Set-TextValue $ws "D12" "0.03343"
Set-TextValue $ws "D13" "0.03054"
Set-TextValue $ws "D14" "0.09336"
Set-TextValue $ws "D15" "3.856"
Set-TextValue $ws "D16" "0.001569"
Set-TextValue $ws "D17" "0.04662"
Set-TextValue $ws "D18" "0.0005953"
Set-TextValue $ws "D19" "0.005886"
Set-TextValue $ws "D20" "0.001265"
Set-TextValue $ws "D21" "0.004899"
Set-TextValue $ws "D22" "0.00006797"
Set-TextValue $ws "D23" "3.566"
Set-TextValue $ws "D24" "2.112"
Set-TextValue $ws "D25" "0.3233"
Set-TextValue $ws "D41" "0.006189"
Set-TextValue $ws "D42" "0.1072"
Set-TextValue $ws "D43" "0.002569"
Set-TextValue $ws "D44" "0.009072"
Set-TextValue $ws "D45" "0.00005218"
Set-TextValue $ws "D47" "0.6705"
Set-TextValue $ws "D48" "0.002385"
